$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update Asset Number (A2), add a Notes/Cost value in F2
$ws.Range("A2").Value = 131166
$ws.Range("F2").Value = 15166

# Row 2: replace the shared strings used for Tune Type (B2) and Staff (E2)
$ws.Range("B2").Value = "Accident Calibration"
$ws.Range("E2").Value = "RG"

# Row 3: clear out the old entry entirely (values only, keep date cell styling)
$ws.Range("A3").ClearContents()
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()

# Update the active selection shown when the workbook is opened
$null = $ws.Range("A2:F4").Select()
